$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = "Interchange Part Number|DOT Compliant Jeep Wrangler "
$ws.Range("P2").Value = "Other Part Number|JW Speaker 8700 Evolution "
$ws.Range("W2").Value = "Part Brand|12V DOT LED High & Low Beam Headlights "
